$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values for rows 2-6
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = -1
